$d = $word.ActiveDocument

# --- Change 1: fix double space after "16." -> single space ---
$d.Content.Find.Execute("16.  Thediagram", $false, $false, $false, $false, $false, $true, 1, $false, "16. Thediagram", 2)

# --- Change 2: replace the inline picture with the OCR'd text lines ---
$ishp = $d.InlineShapes.Item(1)
$p = $ishp.Range.Paragraphs.Item(1)
$r = $p.Range
$ishp.Delete()
$r.InsertBefore(", Support`v`v«— 15cm—rlle— 150m —>`v`vPn nn`v`v|`v`v| ; .`v`v| wooden rod 6 iron ball Q`v`v5`v| magnet P (609)`v| (50g)")

# --- Change 3: garble "Which one of ine following statements Is faise’?" ---
$d.Content.Find.Execute("Which one of ine following statements Is faise’?", $false, $false, $false, $false, $false, $true, 1, $false, "Wnich one oF ine tollowing statements Is faise’?", 2)

# --- Change 4: prefix ": " and add an extra manual line break before "wooden rod will tit" ---
$d.Content.Find.Execute("(4) When the south pole of another magnet is placed. under tron-ball-Q, the", $false, $false, $false, $false, $false, $true, 1, $false, ": (4) When the south pole of another magnet is placed. under tron-ball-Q, the^l", 2)

# --- Change 5: clear the trailing stray-text paragraph, leaving it empty ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$tailRange.Delete()
